$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each crypto row.
# Leading apostrophe forces text so numeric-looking values (e.g. "9.387",
# "152.03", "0.000008069") keep their exact display instead of being
# parsed into floating point numbers.
$ws.Range("D2").Value = "'31.151.05"
$ws.Range("E2").Value = "'  +3.64%  "
$ws.Range("D3").Value = "'1.918.21"
$ws.Range("E3").Value = "'  +1.92%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "'  +0.29%  "
$ws.Range("D5").Value = "'246.04"
$ws.Range("E5").Value = "'  +0.82%  "
$ws.Range("D6").Value = "'0.9988"
$ws.Range("E6").Value = "'  +0.18%  "
$ws.Range("D7").Value = "'0.5004"
$ws.Range("E7").Value = "'  +1.22%  "
$ws.Range("D8").Value = "'0.3013"
$ws.Range("E8").Value = "'  +3.21%  "
$ws.Range("D9").Value = "'0.06910"
$ws.Range("E9").Value = "'  +4.26%  "
$ws.Range("D10").Value = "'1.914.23"
$ws.Range("E10").Value = "'  +1.73%  "
$ws.Range("D11").Value = "'17.06"
$ws.Range("E11").Value = "'  +0.18%  "
$ws.Range("D12").Value = "'0.07308"
$ws.Range("E12").Value = "'  +1.53%  "
$ws.Range("D13").Value = "'90.43"
$ws.Range("E13").Value = "'  +5.65%  "
$ws.Range("D14").Value = "'0.6858"
$ws.Range("E14").Value = "'  +3.15%  "
$ws.Range("D15").Value = "'5.103"
$ws.Range("E15").Value = "'  +4.91%  "
$ws.Range("D16").Value = "'31.101.40"
$ws.Range("E16").Value = "'  +3.55%  "
$ws.Range("D17").Value = "'0.000008069"
$ws.Range("E17").Value = "'  +2.44%  "
$ws.Range("D18").Value = "'13.45"
$ws.Range("E18").Value = "'  +4.95%  "
$ws.Range("D19").Value = "'0.9995"
$ws.Range("E19").Value = "'  +0.24%  "
$ws.Range("D20").Value = "'2.160.99"
$ws.Range("E20").Value = "'  +1.97%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "'  +0.56%  "
$ws.Range("D22").Value = "'4.895"
$ws.Range("E22").Value = "'  +2.65%  "
$ws.Range("D23").Value = "'182.35"
$ws.Range("E23").Value = "'  +34.35%  "
$ws.Range("D24").Value = "'6.130"
$ws.Range("E24").Value = "'  +9.46%  "
$ws.Range("D25").Value = "'9.387"
$ws.Range("E25").Value = "'  +2.52%  "
$ws.Range("D26").Value = "'152.03"
$ws.Range("D27").Value = "'18.25"
$ws.Range("E27").Value = "'  +8.69%  "
$ws.Range("E28").Value = "'  +2.45%  "
$ws.Range("D29").Value = "'1.407"
$ws.Range("E29").Value = "'  +1.86%  "
$ws.Range("D30").Value = "'4.372"
$ws.Range("E30").Value = "'  +4.09%  "
$ws.Range("D31").Value = "'0.08993"
$ws.Range("E31").Value = "'  +3.97%  "
$ws.Range("D32").Value = "'4.074"
$ws.Range("E32").Value = "'  +2.94%  "
$ws.Range("D33").Value = "'0.05264"
$ws.Range("E33").Value = "'  +5.58%  "
$ws.Range("D34").Value = "'0.7559"
$ws.Range("E34").Value = "'  +7.65%  "
$ws.Range("D35").Value = "'1.148"
$ws.Range("E35").Value = "'  +3.26%  "
$ws.Range("D36").Value = "'2.664"
$ws.Range("E36").Value = "'  +0.44%  "
$ws.Range("D37").Value = "'0.01923"
$ws.Range("E37").Value = "'  +16.99%  "
$ws.Range("D38").Value = "'2.752"
$ws.Range("E38").Value = "'  +2.13%  "
$ws.Range("D39").Value = "'2.200"
$ws.Range("E39").Value = "'  -0.55%  "
$ws.Range("D40").Value = "'0.9437"
$ws.Range("E40").Value = "'  +1.11%  "
$ws.Range("D41").Value = "'0.4381"
$ws.Range("E41").Value = "'  +4.12%  "
$ws.Range("D42").Value = "'5.945"
$ws.Range("E42").Value = "'  -0.77%  "
$ws.Range("D43").Value = "'105.38"
$ws.Range("E43").Value = "'  +3.45%  "
$ws.Range("D44").Value = "'7.896"
$ws.Range("E44").Value = "'  +4.03%  "
$ws.Range("D45").Value = "'0.9990"
$ws.Range("E45").Value = "'  +0.12%  "
$ws.Range("D46").Value = "'0.1341"
$ws.Range("E46").Value = "'  +6.29%  "
$ws.Range("D47").Value = "'0.05873"
$ws.Range("E47").Value = "'  +2.83%  "
$ws.Range("D48").Value = "'8.635"
$ws.Range("E48").Value = "'  +4.62%  "
$ws.Range("E49").Value = "'  +5.18%  "
$ws.Range("E50").Value = "'  +2.74%  "
$ws.Range("D51").Value = "'1.393"
$ws.Range("E51").Value = "'  +4.22%  "
